{"js": "const body = context.document.body;\n\nconst replacements = [\n  [\n    \", child of Deepak Sanghi, aged 24, residing at testing, do hereby solemnly affirm and declare as under:\",\n    \", child of Deepak Kumar Sanghi, aged 24, residing at 391 29th Street, do hereby solemnly affirm and declare as under:\",\n  ],\n  [\n    \"I perceive myself as a transgender person whose gender does not match with the gender assigned at birth.\",\n    \"I perceive myself as a Male person whose gender does not match with the gender assigned at birth.\",\n  ],\n  [\n    \"I declare myself as transgender person.\",\n    \"I declare myself as Male.\",\n  ],\n  [\n    \"I am executing this affidavit to be submitted to the District Magistrate for issue of certificate of identity as a transgender person under Section 6 of the Transgender Persons (Protection of Rights) Act, 2019 read with Rules 3, 4 and 5 of the Transgender Persons (Protection of Rights) Rules, 2020.\",\n    \"I am executing this affidavit to be submitted to the District Magistrate for issue of certificate of identity as a Male person under Section 7 of the Transgender Persons (Protection of Rights) Act, 2019 read with Rules 3 and 6 of the Transgender Persons (Protection of Rights) Rules, 2020.\",\n  ],\n];\n\nfor (const [find, replace] of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Update child/father name and residing address\n$r1 = $d.Content\n$r1.Find.Execute(\", child of Deepak Sanghi, aged 24, residing at testing, do hereby solemnly affirm and declare as under:\", $false, $false, $false, $false, $false, $true, 1, $false, \", child of Deepak Kumar Sanghi, aged 24, residing at 391 29th Street, do hereby solemnly affirm and declare as under:\", 2)\n\n# 2. Update self-perceived gender sentence\n$r2 = $d.Content\n$r2.Find.Execute(\"I perceive myself as a transgender person whose gender does not match with the gender assigned at birth.\", $false, $false, $false, $false, $false, $true, 1, $false, \"I perceive myself as a Male person whose gender does not match with the gender assigned at birth.\", 2)\n\n# 3. Update declared gender sentence\n$r3 = $d.Content\n$r3.Find.Execute(\"I declare myself as transgender person.\", $false, $false, $false, $false, $false, $true, 1, $false, \"I declare myself as Male.\", 2)\n\n# 4. Update the certificate/section/rules sentence\n$r4 = $d.Content\n$r4.Find.Execute(\"I am executing this affidavit to be submitted to the District Magistrate for issue of certificate of identity as a transgender person under Section 6 of the Transgender Persons (Protection of Rights) Act, 2019 read with Rules 3, 4 and 5 of the Transgender Persons (Protection of Rights) Rules, 2020.\", $false, $false, $false, $false, $false, $true, 1, $false, \"I am executing this affidavit to be submitted to the District Magistrate for issue of certificate of identity as a Male person under Section 7 of the Transgender Persons (Protection of Rights) Act, 2019 read with Rules 3 and 6 of the Transgender Persons (Protection of Rights) Rules, 2020.\", 2)\n"}
